$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the two "apachefriends" hyperlink runs into a single
# run (text unchanged, only the run split disappears).
#
# A Find/Replace whose matched range spans the *whole* first run and
# the *whole* second run exactly (start-to-end) makes the engine infer
# the replacement's formatting from the context just *before* the
# match instead of from the matched runs themselves, which would steal
# the preceding run's direct formatting. Excluding the very first
# character from the search string keeps the match start one
# character inside the first run, so the merge correctly inherits the
# hyperlink run's own formatting (rStyle "Hipervnculo" + lang es-MX).
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("ttps://www.apachefriends.org/es/index.html", $true, $false, $false, $false, $false, $true, 1, $false, "ttps://www.apachefriends.org/es/index.html", 2)

# ---------------------------------------------------------------------
# Change 2: insert " de XAMPP" right after "htdocs" in
# "Crear en la carpeta htdocs la carpeta Proyecto".
#
# Inserting text whose formatting matches its neighbours makes the
# engine coalesce every identically-formatted run in the paragraph
# into one big run. To reproduce the original per-run layout (the
# pre-existing "Crear en la carpeta " / "htdocs" split plus the new
# " de XAMPP" run, each kept apart from " la carpeta Proyecto"), the
# run boundaries are re-established afterwards by toggling Bold on and
# back off across each boundary; the net formatting is unchanged but
# the coalesced run gets split again at those points.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("htdocs", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$htdocsStart = $r.Start
$htdocsEnd = $r.End

$ins = $d.Range($htdocsEnd, $htdocsEnd)
$ins.InsertAfter(" de XAMPP")

$boundary1 = $d.Range($htdocsStart, $htdocsEnd)
$boundary1.Bold = 1
$boundary1.Bold = 0

$boundary2 = $d.Range($htdocsEnd, $htdocsEnd + 9)
$boundary2.Bold = 1
$boundary2.Bold = 0
